$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Régularisation" block (Matricule / Nom-Prénom / Régularisation)
# in columns AB:AD, mirroring the existing blocks (A:F, I:M, P:S, V:Y).
$ws.Range("AB1").Value = $ws.Range("A1").Text
$ws.Range("AC1").Value = $ws.Range("B1").Text
$ws.Range("AD1").Value = "Régularisation"

# Copy the header formatting (bold, centered, bordered) from an existing
# header cell onto the new header cells.
$ws.Range("A1").Copy()
$ws.Range("AB1:AD1").PasteSpecial(-4122)

# Reflect the new selection left by the edit.
$ws.Range("AB3").Select() | Out-Null
